$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New columns G (7) and H (8) sizing (approximate best-fit widths) ---
$ws.Columns.Item(7).ColumnWidth = 19.1667
$ws.Columns.Item(8).ColumnWidth = 10.5

# --- Row 6 gets an explicit (custom) row height matching the default ---
$ws.Rows.Item(6).RowHeight = 14.6

# --- Copy existing formatting onto the new G:K cells (reuses the workbook's
#     existing style records instead of minting new duplicate ones) ---
$ws.Range("A13").Copy()
$ws.Range("G6:K6").PasteSpecial(-4122)
$ws.Range("A14").Copy()
$ws.Range("G7:K7").PasteSpecial(-4122)
$ws.Range("A9").Copy()
$ws.Range("G9").PasteSpecial(-4122)
$ws.Range("A10").Copy()
$ws.Range("G10").PasteSpecial(-4122)
$ws.Range("A11").Copy()
$ws.Range("G11").PasteSpecial(-4122)
$ws.Range("B9").Copy()
$ws.Range("H9:K9").PasteSpecial(-4122)
$ws.Range("B10").Copy()
$ws.Range("H10:M10").PasteSpecial(-4122)
$ws.Range("B11").Copy()
$ws.Range("H11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- New header text above the "direct construction" block (merged G6:K7) ---
$ws.Range("G6").Value2 = "Here we create Person objects by directly constructing Python classes, no factory function is required."
$ws.Range("G6:K7").Merge()

# --- Row 9: directly-constructed Person objects (G9 label + H9:K9 formulas) ---
$ws.Range("G9").Value2 = "Cached Person Objects"
$ws.Range("H9").FormulaArray = "=_xll.Person(B6,B7,B8)"
$ws.Range("I9").FormulaArray = "=_xll.Person(C6,C7,C8)"
$ws.Range("J9").FormulaArray = "=_xll.Person(D6,D7,D8)"
$ws.Range("K9").FormulaArray = "=_xll.Person(E6,E7,E8)"

# --- Row 10: descriptions of the directly-constructed Person objects ---
$ws.Range("G10").Value2 = "Descriptions"
$ws.Range("H10").FormulaArray = "=_xll.describe(H9)"
$ws.Range("I10").FormulaArray = "=_xll.describe(I9)"
$ws.Range("J10").FormulaArray = "=_xll.describe(J9)"
$ws.Range("K10").FormulaArray = "=_xll.describe(K9)"

# --- Row 11: average age of the directly-constructed Person objects ---
$ws.Range("G11").Value2 = "Average age"
$ws.Range("H11").FormulaArray = "=_xll.average_age(H9:K9)"

# --- View state: update the selected cell shown when the workbook is opened ---
$ws.Range("H3").Select()
